$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35/36: Dai <-> Hedera swap (name, link, price, volume all change)
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.102"
$ws.Range("E35").Value = "  -3.15%  "

$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.13%  "

# Price (D) and Volume(1h) (E) updates for remaining rows
# Row 2
$ws.Range("D2").Value = "66.910.04"
$ws.Range("E2").Value = "  +1.25%  "
# Row 3
$ws.Range("D3").Value = "3.271.31"
$ws.Range("E3").Value = "  -1.90%  "
# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.12%  "
# Row 5
$ws.Range("D5").Value = "577.86"
$ws.Range("E5").Value = "  -1.08%  "
# Row 6
$ws.Range("D6").Value = "172.08"
$ws.Range("E6").Value = "  -6.91%  "
# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.05%  "
# Row 8
$ws.Range("E8").Value = "  +0.50%  "
# Row 9
$ws.Range("D9").Value = "3.270.82"
$ws.Range("E9").Value = "  -1.84%  "
# Row 10
$ws.Range("E10").Value = "  -4.86%  "
# Row 11
$ws.Range("D11").Value = "'0.570"
$ws.Range("E11").Value = "  -1.72%  "
# Row 12
$ws.Range("D12").Value = "44.67"
$ws.Range("E12").Value = "  -4.72%  "
# Row 13
$ws.Range("D13").Value = "'0.0000268"
$ws.Range("E13").Value = "  +0.14%  "
# Row 14
$ws.Range("D14").Value = "681.02"
$ws.Range("E14").Value = "  +3.42%  "
# Row 15
$ws.Range("D15").Value = "3.797.91"
$ws.Range("E15").Value = "  -1.92%  "
# Row 16
$ws.Range("E16").Value = "  -2.93%  "
# Row 17
$ws.Range("D17").Value = "67.046.55"
$ws.Range("E17").Value = "  +1.09%  "
# Row 18
$ws.Range("E18").Value = "  +0.39%  "
# Row 19
$ws.Range("D19").Value = "3.269.88"
$ws.Range("E19").Value = "  -2.08%  "
# Row 20
$ws.Range("D20").Value = "17.15"
$ws.Range("E20").Value = "  -3.84%  "
# Row 21
$ws.Range("D21").Value = "'10.60"
$ws.Range("E21").Value = "  -4.30%  "
# Row 22
$ws.Range("D22").Value = "'0.880"
$ws.Range("E22").Value = "  -1.75%  "
# Row 23
$ws.Range("D23").Value = "16.82"
$ws.Range("E23").Value = "  -4.44%  "
# Row 24
$ws.Range("D24").Value = "5.23"
$ws.Range("E24").Value = "  +3.71%  "
# Row 25
$ws.Range("D25").Value = "98.24"
$ws.Range("E25").Value = "  -2.13%  "
# Row 26
$ws.Range("E26").Value = "  -4.51%  "
# Row 27
$ws.Range("D27").Value = "2.63"
$ws.Range("E27").Value = "  -5.08%  "
# Row 28
$ws.Range("D28").Value = "32.92"
$ws.Range("E28").Value = "  +2.54%  "
# Row 29
$ws.Range("D29").Value = "8.97"
$ws.Range("E29").Value = "  -5.06%  "
# Row 30
$ws.Range("D30").Value = "8.26"
$ws.Range("E30").Value = "  -2.46%  "
# Row 31
$ws.Range("D31").Value = "6.59"
$ws.Range("E31").Value = "  -1.82%  "
# Row 32
$ws.Range("D32").Value = "576.83"
$ws.Range("E32").Value = "  -4.77%  "
# Row 33
$ws.Range("D33").Value = "10.81"
$ws.Range("E33").Value = "  -2.25%  "
# Row 34
$ws.Range("D34").Value = "3.796.17"
$ws.Range("E34").Value = "  -2.29%  "
# Row 37
$ws.Range("D37").Value = "55.36"
$ws.Range("E37").Value = "  -1.62%  "
# Row 38
$ws.Range("D38").Value = "3.27"
$ws.Range("E38").Value = "  -15.20%  "
# Row 39
$ws.Range("D39").Value = "0.128"
$ws.Range("E39").Value = "  -0.85%  "
# Row 40
$ws.Range("D40").Value = "'3.40"
$ws.Range("E40").Value = "  +0.02%  "
# Row 41
$ws.Range("D41").Value = "2.55"
$ws.Range("E41").Value = "  -4.61%  "
# Row 42
$ws.Range("D42").Value = "31.25"
$ws.Range("E42").Value = "  -4.79%  "
# Row 43
$ws.Range("E43").Value = "  -6.44%  "
# Row 44
$ws.Range("D44").Value = "0.323"
$ws.Range("E44").Value = "  -3.73%  "
# Row 45
$ws.Range("D45").Value = "2.96"
$ws.Range("E45").Value = "  -6.79%  "
# Row 46
$ws.Range("D46").Value = "0.0399"
$ws.Range("E46").Value = "  -4.19%  "
# Row 47
$ws.Range("E47").Value = "  -0.25%  "
# Row 48
$ws.Range("E48").Value = "  -1.19%  "
# Row 49
$ws.Range("E49").Value = "  -1.51%  "
# Row 50
$ws.Range("E50").Value = "  -0.22%  "
# Row 51
$ws.Range("D51").Value = "127.56"
$ws.Range("E51").Value = "  -1.24%  "
